$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell C10 value from 18 to 1 (as shown in the commit diff)
$ws.Range("C10").Value = 1
